$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (interested count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 46
$wsExpo.Range("F5").Value = 136
$wsExpo.Range("F6").Value = 9170
$wsExpo.Range("F10").Value = 1059
$wsExpo.Range("F15").Value = 354
$wsExpo.Range("F16").Value = 77
$wsExpo.Range("F17").Value = 245
$wsExpo.Range("F18").Value = 1189

# Sheet "全部类型" (all types) - same events, mirrored rows, update column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 46
$wsAll.Range("F7").Value = 136
$wsAll.Range("F8").Value = 9170
$wsAll.Range("F12").Value = 1059
$wsAll.Range("F17").Value = 354
$wsAll.Range("F18").Value = 77
$wsAll.Range("F19").Value = 245
$wsAll.Range("F20").Value = 1189
